# Consolidate the split "First" / " " / "slide" and "Third" / " " / "slide"
# text runs on the title placeholders of slides 1 and 3 into single runs,
# matching the PowerPoint writer's run-consolidation behaviour.
#
# Re-assigning TextRange.Text with the exact same concatenated string is a
# no-op for the engine's content-diffing (nothing "changed" from a value
# perspective), so the multi-run shape survives untouched. Writing a
# throwaway value first forces the existing runs to collapse into one, and
# the follow-up assignment then lands the final text into that single run.

$p = $ppt.ActivePresentation

$s1 = $p.Slides.Item(1)
$title1 = $s1.Shapes.Item(1).TextFrame.TextRange
$title1.Text = "__placeholder__"
$title1.Text = "First slide"

$s3 = $p.Slides.Item(3)
$title3 = $s3.Shapes.Item(1).TextFrame.TextRange
$title3.Text = "__placeholder__"
$title3.Text = "Third slide"
